$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "nan"/empty exit velo + launch angle cells
$ws.Range("M10").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("M30").Value = ""
$ws.Range("M37").Value = ""
$ws.Range("M39").Value = ""

# Pitch mix changes
$ws.Range("J17").Value = "FB,CB,CH"
$ws.Range("J26").Value = "FB,CB,CH"
$ws.Range("J35").Value = "FB,CB,CH"
$ws.Range("J44").Value = "SL,FB,CB,CH"

# Inning numbers
$ws.Range("J19").Value = 3
$ws.Range("J28").Value = 4
$ws.Range("J37").Value = 5

# Outs
$ws.Range("J29").Value = 1
$ws.Range("J38").Value = 2

# Hit type / result
$ws.Range("M23").Value = "Fly Ball"
$ws.Range("M24").Value = "Out"
$ws.Range("M32").Value = "Line Drive"
$ws.Range("M33").Value = "Single"
$ws.Range("M41").Value = "Ground Ball"
$ws.Range("M42").Value = "Double"

# Pitcher names
$ws.Range("J32").Value = "Roblez"
$ws.Range("J41").Value = "Herbst"

# FB Velo
$ws.Range("J34").Value = "88-90 MPH"
$ws.Range("J43").Value = "83-85 MPH"
